$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.547.76"
$ws.Range("E2").Value = "  +4.71%  "

$ws.Range("D3").Value = "3.498.34"
$ws.Range("E3").Value = "  +2.84%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.10"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.88%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.86"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +7.44%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").Value = "3.499.32"
$ws.Range("E8").Value = "  +2.85%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.579"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.84%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.26"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.10%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.124"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +4.15%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.435"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.31%  "

$ws.Range("D13").Value = "4.103.44"
$ws.Range("E13").Value = "  +2.87%  "

$ws.Range("E14").Value = "  +0.46%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.92"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +3.46%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "66.498.49"
$ws.Range("E16").Value = "  +4.49%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000177"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +3.12%  "

$ws.Range("D18").Value = "3.506.99"
$ws.Range("E18").Value = "  +2.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.27"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.93%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.99"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.29%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "387.99"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.82%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.95"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.91%  "

$ws.Range("E23").Value = "  +2.55%  "

$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.526"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.82%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000123"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +7.68%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.12"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +4.42%  "

$ws.Range("E28").Value = "  +1.75%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.03%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.33"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +4.60%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.46"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +5.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.05"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.79%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.41"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.23%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.39"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +6.64%  "

$ws.Range("E35").Value = "  +0.07%  "

$ws.Range("E36").Value = "  +3.26%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "161.26"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.31%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.901"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +7.45%  "

$ws.Range("E39").Value = "  +5.59%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0746"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.74%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.74"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +5.32%  "

$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.64"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +6.30%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.30"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.96%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.86"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +4.90%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "43.44"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.30%  "

$ws.Range("D46").Value = "2.785.42"
$ws.Range("E46").Value = "  -1.02%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0313"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.95%  "

$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.52"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +8.63%  "

$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "351.90"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +7.83%  "

$ws.Range("E50").Value = "  +5.85%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.25"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +10.92%  "
